$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> values for columns D (Fecha serial), L (Calidad), M (Volumen),
# N (Precio minimo), O (Precio maximo), P (Precio promedio ponderado),
# Q (Unidad de comercializacion), S (Precio $/Kg), T (Kg / unidad)
$rows = @{
    2  = @{ D = 44553; L = "Especial"; M = 200; N = 22000; O = 22000; P = 22000; Q = "`$/bandeja 6 kilos"; S = 3667; T = 6 }
    3  = @{ D = 44553; L = "Primera";  M = 150; N = 18000; O = 18000; P = 18000; Q = "`$/bandeja 6 kilos"; S = 3000; T = 6 }
    4  = @{ D = 44189; L = "Especial"; M = 20;  N = 15000; O = 15000; P = 15000; Q = "`$/bandeja 7 kilos"; S = 2143; T = 7 }
    5  = @{ D = 44189; L = "Primera";  M = 30;  N = 13000; O = 13000; P = 13000; Q = "`$/bandeja 7 kilos"; S = 1857; T = 7 }
    6  = @{ D = 44550; L = "Primera";  M = 60;  N = 24000; O = 24000; P = 24000; Q = "`$/bandeja 7 kilos"; S = 3429; T = 7 }
    7  = @{ D = 44561; L = "Primera";  M = 200; N = 18000; O = 18000; P = 18000; Q = "`$/bandeja 6 kilos"; S = 3000; T = 6 }
    8  = @{ D = 44204; L = "Primera";  M = 110; N = 7000;  O = 7500;  P = 7318;  Q = "`$/bandeja 7 kilos"; S = 1045; T = 7 }
    9  = @{ D = 44572; L = "Primera";  M = 65;  N = 20000; O = 20000; P = 20000; Q = "`$/bandeja 6 kilos"; S = 3333; T = 6 }
    10 = @{ D = 44187; L = "Especial"; M = 45;  N = 14000; O = 14000; P = 14000; Q = "`$/bandeja 7 kilos"; S = 2000; T = 7 }
    11 = @{ D = 44187; L = "Primera";  M = 50;  N = 12000; O = 12000; P = 12000; Q = "`$/bandeja 7 kilos"; S = 1714; T = 7 }
    12 = @{ D = 44558; L = "Especial"; M = 20;  N = 22000; O = 22000; P = 22000; Q = "`$/bandeja 6 kilos"; S = 3667; T = 6 }
    13 = @{ D = 44558; L = "Primera";  M = 25;  N = 18000; O = 18000; P = 18000; Q = "`$/bandeja 6 kilos"; S = 3000; T = 6 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 4).Value2  = $vals.D   # D: Fecha
    $ws.Cells.Item($r, 12).Value  = $vals.L   # L: Calidad
    $ws.Cells.Item($r, 13).Value  = $vals.M   # M: Volumen
    $ws.Cells.Item($r, 14).Value  = $vals.N   # N: Precio minimo
    $ws.Cells.Item($r, 15).Value  = $vals.O   # O: Precio maximo
    $ws.Cells.Item($r, 16).Value  = $vals.P   # P: Precio promedio ponderado
    $ws.Cells.Item($r, 17).Value  = $vals.Q   # Q: Unidad de comercializacion
    $ws.Cells.Item($r, 19).Value  = $vals.S   # S: Precio $/Kg
    $ws.Cells.Item($r, 20).Value  = $vals.T   # T: Kg / unidad
}
